$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 rework ---------------------------------------------------------
# Build a throw-away text cell holding the literal string "1" (as a real
# text value, not a number) by wrapping it in T(), then copy *values only*
# from it into B4 and H4 so both become shared-string cells containing "1"
# (matches the diff: B4/H4 -> t="s" pointing at the "1" string) without
# dragging any formatting/styles along.
$ws.Range("Z1").Formula = '=T("1")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4163) | Out-Null
$ws.Range("H4").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").ClearContents() | Out-Null

# The rest of the old "DataV10" block (C4:F4) and the old Trama_Completa
# trio (I4:K4) lose their values; I4:K4 keep their existing style (s="4").
$ws.Range("C4:F4").ClearContents() | Out-Null
$ws.Range("I4:K4").ClearContents() | Out-Null

# --- Column B width --------------------------------------------------------
# Drop the bestFit auto-width in favour of an explicit custom width.
$ws.Columns("B").ColumnWidth = 20

# --- Selection / view -------------------------------------------------------
$ws.Range("H4").Select() | Out-Null
